$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Novo loyout das curtidas" — fill in the weekly sprint ("Sprint semanal")
# column J for the rows that were missing it, and re-style the remaining
# blank row J24 (underlined font) — "sistema sanduiche dos comentarios".

$ws.Range("J3").Value = "Semana 3"
$ws.Range("J6").Value = "Semana 2"
$ws.Range("J7").Value = "Semana 3"
$ws.Range("J8").Value = "Semana 3"
$ws.Range("J9").Value = "Semana 3"
$ws.Range("J10").Value = "Semana 2"
$ws.Range("J11").Value = "Semana 3"
$ws.Range("J12").Value = "Semana 4"
$ws.Range("J18").Value = "Semana 4"
$ws.Range("J19").Value = "Semana 4"
$ws.Range("J20").Value = "Semana 3"
$ws.Range("J21").Value = "Semana 3"
$ws.Range("J22").Value = "Semana 3"
$ws.Range("J23").Value = "Semana 3"
$ws.Range("J25").Value = "Semana 3"

# J24 stays blank but gets a new underlined-font style.
$ws.Range("J24").Font.Underline = $true

# Recalculate so the SUMIF totals in M5:M7 pick up the new weeks.
$excel.Calculate()

# Move the selection to where the author ended up.
[void]$ws.Range("H64").Select()

Write-Host "done"
